$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 167; existing rows 167:199 shift down to 168:200
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new data record
$ws.Cells.Item(167, 1).Value = 3
$ws.Cells.Item(167, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(167, 3).Value = "Coquimbo"
$ws.Cells.Item(167, 4).Value = 44476
$ws.Cells.Item(167, 5).Value = 5
$ws.Cells.Item(167, 6).Value = 100114013
$ws.Cells.Item(167, 7).Value = "Zanahoria"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 230
$ws.Cells.Item(167, 11).Value = 8000
$ws.Cells.Item(167, 12).Value = 8500
$ws.Cells.Item(167, 13).Value = 8261
$ws.Cells.Item(167, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(167, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(167, 16).Value = 413
$ws.Cells.Item(167, 17).Value = 20
$ws.Cells.Item(167, 18).Value = "Hortaliza"

# Match the date-style formatting used by column D elsewhere in the sheet
$ws.Cells.Item(167, 4).NumberFormat = $ws.Cells.Item(168, 4).NumberFormat
